$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing daily data runs through row 251 (date serial 44325).
# Append 4 more rows (252-255) continuing the sequence with zero counts,
# matching the existing date-column formatting ("aggiornamento fino a 13/03").

$lastRow = 251
$startSerial = 44326
$newRowCount = 4

$dateTemplate = $ws.Cells.Item($lastRow, 1)

for ($i = 0; $i -lt $newRowCount; $i++) {
    $row = $lastRow + 1 + $i
    $serial = $startSerial + $i

    $a = $ws.Cells.Item($row, 1)
    $a.Value = $serial

    # Copy the date column's style/format (border, bold, alignment,
    # number format) from the previous row instead of recreating it,
    # so no duplicate style entries are introduced.
    $dateTemplate.Copy()
    $a.PasteSpecial(-4122)

    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

$excel.CutCopyMode = $false
